# Add "Positive Instances" and "Negative Instances" columns (F and G) to the
# TableDatasets sheet, shifting the existing Features / Dimensionality /
# Balance columns two places to the right (F->H, G->I, H->J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, blank columns at F:G. This shifts the old F, G, H columns
# (Features, Dimensionality, Balance) to H, I, J respectively, and carries
# the header cell style (bold/centered/bordered) along with it.
$ws.Columns("F:G").Insert()

# Fill in the new header labels.
$ws.Range("F1").Value = "Positive Instances"
$ws.Range("G1").Value = "Negative Instances"

# Per-dataset positive / negative instance counts (rows 2-51).
$counts = @(
    @{Row=2; Pos=63; Neg=51},
    @{Row=3; Pos=111; Neg=57},
    @{Row=4; Pos=301; Neg=276},
    @{Row=5; Pos=35; Neg=23},
    @{Row=6; Pos=46; Neg=24},
    @{Row=7; Pos=76; Neg=14},
    @{Row=8; Pos=31; Neg=88},
    @{Row=9; Pos=94; Neg=167},
    @{Row=10; Pos=136; Neg=143},
    @{Row=11; Pos=33; Neg=21},
    @{Row=12; Pos=41; Neg=47},
    @{Row=13; Pos=38; Neg=13},
    @{Row=14; Pos=12; Neg=72},
    @{Row=15; Pos=25; Neg=68},
    @{Row=16; Pos=61; Neg=30},
    @{Row=17; Pos=61; Neg=76},
    @{Row=18; Pos=159; Neg=134},
    @{Row=19; Pos=60; Neg=147},
    @{Row=20; Pos=133; Neg=50},
    @{Row=21; Pos=97; Neg=115},
    @{Row=22; Pos=282; Neg=238},
    @{Row=23; Pos=92; Neg=69},
    @{Row=24; Pos=119; Neg=154},
    @{Row=25; Pos=102; Neg=57},
    @{Row=26; Pos=115; Neg=58},
    @{Row=27; Pos=32; Neg=19},
    @{Row=28; Pos=55; Neg=20},
    @{Row=29; Pos=38; Neg=50},
    @{Row=30; Pos=23; Neg=121},
    @{Row=31; Pos=32; Neg=35},
    @{Row=32; Pos=637; Neg=332},
    @{Row=33; Pos=47; Neg=81},
    @{Row=34; Pos=596; Neg=177},
    @{Row=35; Pos=44; Neg=15},
    @{Row=36; Pos=50; Neg=41},
    @{Row=37; Pos=68; Neg=70},
    @{Row=38; Pos=127; Neg=133},
    @{Row=39; Pos=373; Neg=45},
    @{Row=40; Pos=79; Neg=108},
    @{Row=41; Pos=47; Neg=103},
    @{Row=42; Pos=37; Neg=40},
    @{Row=43; Pos=72; Neg=131},
    @{Row=44; Pos=124; Neg=121},
    @{Row=45; Pos=57; Neg=57},
    @{Row=46; Pos=94; Neg=92},
    @{Row=47; Pos=47; Neg=48},
    @{Row=48; Pos=27; Neg=40},
    @{Row=49; Pos=104; Neg=99},
    @{Row=50; Pos=145; Neg=110},
    @{Row=51; Pos=126; Neg=66}
)

foreach ($entry in $counts) {
    $r = $entry.Row
    $ws.Cells.Item($r, 6).Value = $entry.Pos
    $ws.Cells.Item($r, 7).Value = $entry.Neg
}
